$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2 through 421 from 45205 to 45206
$ws.Range("C2:C421").Value = 45206
